$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint_One")

# Increase the sprint's start value from 71 to 77 (+6, for three new tasks).
$ws.Range("B3").Value = 77

# Move the active selection to B4, matching the saved cursor position.
$ws.Range("B4").Select()
